$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreation")

$ws.Range("A2").Value = "Marina"
$ws.Range("B2").Value = "Avery"
$ws.Range("C2").Value = "testaccount25@email.com"
$ws.Range("D2").Value = "testaccount25"
$ws.Range("J2").Value = "'8173678441"

$ws.Range("J2").Select()
